$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the three new worksheets, in order, at the end of the workbook:
#    CypherOutput_Message, StatOutput, StatOutput_Message
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$wsCypherMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCypherMsg.Name = "CypherOutput_Message"

$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$wsStat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsStat.Name = "StatOutput"

$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$wsStatMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsStatMsg.Name = "StatOutput_Message"

# ---------------------------------------------------------------------------
# 2. CypherOutput_Message: same 10-row "connection log" as the Message sheet
# ---------------------------------------------------------------------------
$wsCypherMsg.Range("A1").Value = 'Neo4j_URL:'
$wsCypherMsg.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$wsCypherMsg.Range("A3").Value = 'User_name:'
$wsCypherMsg.Range("A4").Value = 'neo4j'
$wsCypherMsg.Range("A5").Value = 'PWD:'
$wsCypherMsg.Range("A6").Value = 'icdcDBneo4j0'
$wsCypherMsg.Range("A7").Value = 'Cypher:'
$wsCypherMsg.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN [''Samoyed''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$wsCypherMsg.Range("A9").Value = 'Output:'
$wsCypherMsg.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC35_Canine_Filter_Breed-Samoyed_Neo4jData.xlsx'

# ---------------------------------------------------------------------------
# 3. StatOutput: header row of stat names + a row of counts
#    "0"/"1" must land as *text* (shared-string) cells, not numbers -- a
#    leading apostrophe forces Excel to keep them as literal text, exactly
#    like typing '0 / '1 into a cell.
# ---------------------------------------------------------------------------
$wsStat.Range("A1").Value = 'number_of_files'
$wsStat.Range("B1").Value = 'number_of_sample'
$wsStat.Range("C1").Value = 'number_of_cases'
$wsStat.Range("D1").Value = 'number_of_study'
$wsStat.Range("A2").Value = '''0'
$wsStat.Range("B2").Value = '''0'
$wsStat.Range("C2").Value = '''1'
$wsStat.Range("D2").Value = '''1'

# ---------------------------------------------------------------------------
# 4. StatOutput_Message: connection log (rows 1-7 twice), then the new
#    cypher query (with OPTIONAL MATCH / counts) on row 18, then Output:
#    and the output path.
# ---------------------------------------------------------------------------
$wsStatMsg.Range("A1").Value = 'Neo4j_URL:'
$wsStatMsg.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$wsStatMsg.Range("A3").Value = 'User_name:'
$wsStatMsg.Range("A4").Value = 'neo4j'
$wsStatMsg.Range("A5").Value = 'PWD:'
$wsStatMsg.Range("A6").Value = 'icdcDBneo4j0'
$wsStatMsg.Range("A7").Value = 'Cypher:'
$wsStatMsg.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN [''Samoyed''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$wsStatMsg.Range("A9").Value = 'Output:'
$wsStatMsg.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC35_Canine_Filter_Breed-Samoyed_Neo4jData.xlsx'
$wsStatMsg.Range("A11").Value = 'Neo4j_URL:'
$wsStatMsg.Range("A12").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$wsStatMsg.Range("A13").Value = 'User_name:'
$wsStatMsg.Range("A14").Value = 'neo4j'
$wsStatMsg.Range("A15").Value = 'PWD:'
$wsStatMsg.Range("A16").Value = 'icdcDBneo4j0'
$wsStatMsg.Range("A17").Value = 'Cypher:'
$wsStatMsg.Range("A18").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Samoyed'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$wsStatMsg.Range("A19").Value = 'Output:'
$wsStatMsg.Range("A20").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC35_Canine_Filter_Breed-Samoyed_Neo4jData.xlsx'

# ---------------------------------------------------------------------------
# 5. Restore the original active sheet/tab so the workbook-level view state
#    (tabSelected / activeTab) stays on CypherOutput, like before the edit.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
